$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 51.59157666666666
$ws.Range("H2").Value = 154.77473
$ws.Range("I2").Value = 0.2641250550177587
$ws.Range("J2").Value = 0.2641250550177588
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.441269
$ws.Range("N2").Value = 40.323807
$ws.Range("O2").Value = 0.0897308213348123
$ws.Range("P2").Value = 0.08973082133481232
$ws.Range("Q2").Value = 693.4562601107899
$ws.Range("R2").Value = 6241.106340997109
$ws.Range("S2").Value = 0.02370015812184598
$ws.Range("T2").Value = 0.02370015812184599

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 51.59157666666666
$ws.Range("H3").Value = 154.77473
$ws.Range("I3").Value = 0.2641250550177587
$ws.Range("J3").Value = 0.2641250550177588
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 54.711535
$ws.Range("N3").Value = 164.134605
$ws.Range("O3").Value = 0.3652416280068742
$ws.Range("P3").Value = 0.3652416280068742
$ws.Range("Q3").Value = 2822.654352503516
$ws.Range("R3").Value = 25403.88917253165
$ws.Range("S3").Value = 0.09646946509209141
$ws.Range("T3").Value = 0.09646946509209144

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 51.59157666666666
$ws.Range("H4").Value = 154.77473
$ws.Range("I4").Value = 0.2641250550177587
$ws.Range("J4").Value = 0.2641250550177588
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 63.67711
$ws.Range("N4").Value = 191.03133
$ws.Range("O4").Value = 0.4250937452800914
$ws.Range("P4").Value = 0.4250937452800915
$ws.Range("Q4").Value = 3285.202502476766
$ws.Range("R4").Value = 29566.82252229089
$ws.Range("S4").Value = 0.1122779088598093
$ws.Range("T4").Value = 0.1122779088598093

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 51.59157666666666
$ws.Range("H5").Value = 154.77473
$ws.Range("I5").Value = 0.2641250550177587
$ws.Range("J5").Value = 0.2641250550177588
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.96553866666667
$ws.Range("N5").Value = 53.896616
$ws.Range("O5").Value = 0.119933805378222
$ws.Range("P5").Value = 0.119933805378222
$ws.Range("Q5").Value = 926.8704654792978
$ws.Range("R5").Value = 8341.834189313679
$ws.Range("S5").Value = 0.03167752294401206
$ws.Range("T5").Value = 0.03167752294401206

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09891538535728452
$ws.Range("J6").Value = 0.09891538535728453
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.441269
$ws.Range("N6").Value = 40.323807
$ws.Range("O6").Value = 0.0897308213348123
$ws.Range("P6").Value = 0.08973082133481232
$ws.Range("Q6").Value = 259.7008193435801
$ws.Range("R6").Value = 2337.30737409222
$ws.Range("S6").Value = 0.008875758770758605
$ws.Range("T6").Value = 0.008875758770758609

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09891538535728452
$ws.Range("J7").Value = 0.09891538535728453
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 54.711535
$ws.Range("N7").Value = 164.134605
$ws.Range("O7").Value = 0.3652416280068742
$ws.Range("P7").Value = 0.3652416280068742
$ws.Range("Q7").Value = 1057.089956837033
$ws.Range("R7").Value = 9513.809611533301
$ws.Range("S7").Value = 0.03612801638282192
$ws.Range("T7").Value = 0.03612801638282193

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09891538535728452
$ws.Range("J8").Value = 0.09891538535728453
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 63.67711
$ws.Range("N8").Value = 191.03133
$ws.Range("O8").Value = 0.4250937452800914
$ws.Range("P8").Value = 0.4250937452800915
$ws.Range("Q8").Value = 1230.315206133534
$ws.Range("R8").Value = 11072.8368552018
$ws.Range("S8").Value = 0.04204831162735159
$ws.Range("T8").Value = 0.0420483116273516

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09891538535728452
$ws.Range("J9").Value = 0.09891538535728453
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.96553866666667
$ws.Range("N9").Value = 53.896616
$ws.Range("O9").Value = 0.119933805378222
$ws.Range("P9").Value = 0.119933805378222
$ws.Range("Q9").Value = 347.1149272945956
$ws.Range("R9").Value = 3124.034345651361
$ws.Range("S9").Value = 0.01186329857635239
$ws.Range("T9").Value = 0.01186329857635239

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 112.3724673333333
$ws.Range("H10").Value = 337.117402
$ws.Range("I10").Value = 0.5752951554216499
$ws.Range("J10").Value = 0.57529515542165
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.441269
$ws.Range("N10").Value = 40.323807
$ws.Range("O10").Value = 0.0897308213348123
$ws.Range("P10").Value = 0.08973082133481232
$ws.Range("Q10").Value = 1510.428561621046
$ws.Range("R10").Value = 13593.85705458941
$ws.Range("S10").Value = 0.05162170680592314
$ws.Range("T10").Value = 0.05162170680592316

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 112.3724673333333
$ws.Range("H11").Value = 337.117402
$ws.Range("I11").Value = 0.5752951554216499
$ws.Range("J11").Value = 0.57529515542165
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 54.711535
$ws.Range("N11").Value = 164.134605
$ws.Range("O11").Value = 0.3652416280068742
$ws.Range("P11").Value = 0.3652416280068742
$ws.Range("Q11").Value = 6148.070179544023
$ws.Range("R11").Value = 55332.6316158962
$ws.Range("S11").Value = 0.2101217391506711
$ws.Range("T11").Value = 0.2101217391506712

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 112.3724673333333
$ws.Range("H12").Value = 337.117402
$ws.Range("I12").Value = 0.5752951554216499
$ws.Range("J12").Value = 0.57529515542165
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 63.67711
$ws.Range("N12").Value = 191.03133
$ws.Range("O12").Value = 0.4250937452800914
$ws.Range("P12").Value = 0.4250937452800915
$ws.Range("Q12").Value = 7155.553963356072
$ws.Range("R12").Value = 64399.98567020465
$ws.Range("S12").Value = 0.2445543722596814
$ws.Range("T12").Value = 0.2445543722596815

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 112.3724673333333
$ws.Range("H13").Value = 337.117402
$ws.Range("I13").Value = 0.5752951554216499
$ws.Range("J13").Value = 0.57529515542165
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.96553866666667
$ws.Range("N13").Value = 53.896616
$ws.Range("O13").Value = 0.119933805378222
$ws.Range("P13").Value = 0.119933805378222
$ws.Range("Q13").Value = 2018.831906945737
$ws.Range("R13").Value = 18169.48716251163
$ws.Range("S13").Value = 0.06899733720537414
$ws.Range("T13").Value = 0.06899733720537415

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 12.044915
$ws.Range("H14").Value = 36.134745
$ws.Range("I14").Value = 0.06166440420330686
$ws.Range("J14").Value = 0.06166440420330688
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 13.441269
$ws.Range("N14").Value = 40.323807
$ws.Range("O14").Value = 0.0897308213348123
$ws.Range("P14").Value = 0.08973082133481232
$ws.Range("Q14").Value = 161.898942597135
$ws.Range("R14").Value = 1457.090483374215
$ws.Range("S14").Value = 0.005533197636284577
$ws.Range("T14").Value = 0.005533197636284579

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 12.044915
$ws.Range("H15").Value = 36.134745
$ws.Range("I15").Value = 0.06166440420330686
$ws.Range("J15").Value = 0.06166440420330688
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 54.711535
$ws.Range("N15").Value = 164.134605
$ws.Range("O15").Value = 0.3652416280068742
$ws.Range("P15").Value = 0.3652416280068742
$ws.Range("Q15").Value = 658.995788594525
$ws.Range("R15").Value = 5930.962097350725
$ws.Range("S15").Value = 0.02252240738128973
$ws.Range("T15").Value = 0.02252240738128974

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 12.044915
$ws.Range("H16").Value = 36.134745
$ws.Range("I16").Value = 0.06166440420330686
$ws.Range("J16").Value = 0.06166440420330688
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 63.67711
$ws.Range("N16").Value = 191.03133
$ws.Range("O16").Value = 0.4250937452800914
$ws.Range("P16").Value = 0.4250937452800915
$ws.Range("Q16").Value = 766.9853773956501
$ws.Range("R16").Value = 6902.868396560851
$ws.Range("S16").Value = 0.02621315253324913
$ws.Range("T16").Value = 0.02621315253324914

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 12.044915
$ws.Range("H17").Value = 36.134745
$ws.Range("I17").Value = 0.06166440420330686
$ws.Range("J17").Value = 0.06166440420330688
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.96553866666667
$ws.Range("N17").Value = 53.896616
$ws.Range("O17").Value = 0.119933805378222
$ws.Range("P17").Value = 0.119933805378222
$ws.Range("Q17").Value = 216.3933861692134
$ws.Range("R17").Value = 1947.54047552292
$ws.Range("S17").Value = 0.00739564665248342
$ws.Range("T17").Value = 0.007395646652483423
